$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$model  = $wb.Worksheets.Item("model")

# ---------------------------------------------------------------------------
# Sheet "model": append the first batch of new field definitions
# (gender, age, flag, distributions) right after the existing table.
# ---------------------------------------------------------------------------
$lastDataRow = 23
$formatSrcA = $model.Cells.Item($lastDataRow, 1)
$formatSrcB = $model.Cells.Item($lastDataRow, 2)

function Add-ModelRow([int]$row, [string]$name) {
  $cellA = $model.Cells.Item($row, 1)
  $cellB = $model.Cells.Item($row, 2)

  $cellA.Value = "string"
  $cellB.Value = $name

  $formatSrcA.Copy()
  $cellA.PasteSpecial(-4122)
  $formatSrcB.Copy()
  $cellB.PasteSpecial(-4122)
}

Add-ModelRow 24 "gender"
Add-ModelRow 25 "age"
Add-ModelRow 26 "flag"
Add-ModelRow 27 "distributions"

# ---------------------------------------------------------------------------
# Rename the envelope-barcode question to an item-barcode question:
# "model" B13, then "survey" rows 3 & 4 (field name, prompt, hint).
# ---------------------------------------------------------------------------
$model.Cells.Item(13,2).Value = "item_code"

$survey.Cells.Item(3,2).Value = "item_code"
$survey.Cells.Item(3,3).Value = "Scan their item barcode."
$survey.Cells.Item(4,2).Value = "item_code"
$survey.Cells.Item(4,3).Value = "Update or enter their item code."
$survey.Cells.Item(3,4).Value = "Item code"
$survey.Cells.Item(4,4).Value = "Item code"

# ---------------------------------------------------------------------------
# Sheet "model": append the second batch of new field definitions.
# ---------------------------------------------------------------------------
Add-ModelRow 28 "distributions_received"
Add-ModelRow 29 "remaining_distributions"
Add-ModelRow 30 "scanned_item_code"
Add-ModelRow 31 "mobile_provider"

# Last row of the table: same as above, but closes the box with a top
# border (like the header row does), separating it from the empty rows
# below.
Add-ModelRow 32 "hh_size"
$lastRow = 32
$lastCellA = $model.Cells.Item($lastRow, 1)
$lastCellB = $model.Cells.Item($lastRow, 2)
$lastCellA.Borders.Item(8).LineStyle = 1
$lastCellA.Borders.Item(8).Weight = 2
$lastCellB.Borders.Item(8).LineStyle = 1
$lastCellB.Borders.Item(8).Weight = 2

$model.Activate()
$model.Range("C31").Select()

$excel.CutCopyMode = 0
